$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; this shifts the existing A:E data
# (including the blank "missing data" marker cells) one column to the right, B:F.
$ws.Columns.Item(1).Insert()

# Give the new A1 header cell the same formatting (bold, centered, bordered) as the
# other header cells, then set its text to "ID".
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(1, 1).Value = "ID"

# Populate the new ID column with the row labels.
$ws.Cells.Item(2, 1).Value = "Hb 2"
$ws.Cells.Item(3, 1).Value = "Hb 3"
$ws.Cells.Item(4, 1).Value = "S 24"
$ws.Cells.Item(5, 1).Value = "S 28"
$ws.Cells.Item(6, 1).Value = "Hb 107"
$ws.Cells.Item(7, 1).Value = "Hb 66"
$ws.Cells.Item(8, 1).Value = "Hb 69"
$ws.Cells.Item(9, 1).Value = "Hb 95"
$ws.Cells.Item(10, 1).Value = "Hb 99"
$ws.Cells.Item(11, 1).Value = "Hb 92"
$ws.Cells.Item(12, 1).Value = "Hb 40"
$ws.Cells.Item(13, 1).Value = "Hb 41"
$ws.Cells.Item(14, 1).Value = "S 11"
$ws.Cells.Item(15, 1).Value = "Hb 57"
$ws.Cells.Item(16, 1).Value = "S 21"
$ws.Cells.Item(17, 1).Value = "S 22"
$ws.Cells.Item(18, 1).Value = "S 3"
$ws.Cells.Item(19, 1).Value = "S 4"
$ws.Cells.Item(20, 1).Value = "S 5"
$ws.Cells.Item(21, 1).Value = "Hb 74"
$ws.Cells.Item(22, 1).Value = "Hb 79"
$ws.Cells.Item(23, 1).Value = "Hb 32"
$ws.Cells.Item(24, 1).Value = "S 15"
$ws.Cells.Item(25, 1).Value = "S 16"
